$d = $word.ActiveDocument

# Fix the typo "рессурса" -> "ресурса" (double "с" -> single "с") in the
# answer to question 3 ("... определяет MIME тип рессурса ").
$d.Content.Find.Execute("рессурса", $false, $false, $false, $false, $false,
                         $true, 1, $false, "ресурса", 2)
